# Setting up run modes for test suites
# Adds a new "test_suite" sheet (first tab) to testdata.xlsx with a
# TCID / Runmode table that lets the developer flip Y/N per test suite.

$wb = $excel.ActiveWorkbook

# --- Insert the new "test_suite" sheet as the first tab -------------------
$firstSheet = $wb.Worksheets.Item(1)
$ts = $wb.Worksheets.Add($firstSheet)
$ts.Name = "test_suite"

# --- Populate the TCID / Runmode table -------------------------------------
# Header row first, then the whole Runmode column defaulted to "Y", then the
# TCID names, and finally flip AddCustomer's run mode to "N" -- this mirrors
# how the sheet was actually built and keeps the shared-string insertion
# order lined up with the source workbook.
$ts.Range("A1").Value = "TCID"
$ts.Range("B1").Value = "Runmode"

$ts.Range("B2").Value = "Y"
$ts.Range("B3").Value = "Y"
$ts.Range("B4").Value = "Y"
$ts.Range("B5").Value = "Y"

$ts.Range("A2").Value = "BankManagerLogin"
$ts.Range("A3").Value = "AddCustomer"
$ts.Range("A4").Value = "CustomerLogin"
$ts.Range("A5").Value = "OpenAccount"

$ts.Range("B3").Value = "N"

# Match the page margins used elsewhere in the workbook for freshly added
# sheets (~1.3cm/2cm/0.8cm).
$ts.PageSetup.LeftMargin = 36.850393728
$ts.PageSetup.RightMargin = 36.850393728
$ts.PageSetup.TopMargin = 56.6929134
$ts.PageSetup.BottomMargin = 56.6929134
$ts.PageSetup.HeaderMargin = 22.67716464
$ts.PageSetup.FooterMargin = 22.67716464

# Selection ends on the last populated cell.
$ts.Range("B5").Select()

# --- Selection tweak on the addCustomer sheet ------------------------------
$addCustomer = $wb.Worksheets.Item("addCustomer")
$addCustomer.Range("D1").Select()

# --- Window / view state ----------------------------------------------------
$excel.ActiveWindow.WindowState = -4137  # xlMaximized
$addCustomer.Activate()
